$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Name" column used to hold a short nickname ("BahSh"); the report now
# builds the full display name out of the person's first and last name
# ("Bahronov Shaxriyor") instead.
$ws.Range("B2").Value = "Bahronov Shaxriyor"

# Column B now needs to be wide enough to show the longer name.
$ws.Columns.Item(2).AutoFit()

# The balance column (Saqlash - Avans - Astatka) is now stored as a plain
# computed value instead of a live formula.
$ws.Range("F2").Value = 2131899
$ws.Range("F3").Value = 90
$ws.Range("F4").Value = 5008847

# Leave the selection on F2, matching the last-active cell in the sheet.
$ws.Range("F2").Select()
